$d = $word.ActiveDocument
$sel = $word.Selection

# Move to the very end of the document (after "CSS or flex box course UDEMY")
$sel.EndKey(6)

# A placeholder character is typed into each "empty" paragraph and then
# stripped back out afterwards. Typing a paragraph mark alone (TypeParagraph
# with no TypeText before it) leaves behind a formatted-but-empty <w:r/>
# in the saved XML; going through a placeholder + Range.Text = "" avoids
# that and yields a <w:p> with no run at all, matching a plain blank line.
$placeholder = "#TMP#"

$sel.TypeParagraph()
$sel.TypeText($placeholder)

$sel.TypeParagraph()
$sel.TypeText($placeholder)

$sel.TypeParagraph()
$sel.TypeText("Auxillary function")

$sel.TypeParagraph()
$sel.TypeText("Dynamic css classes with component scope only")

$sel.TypeParagraph()
$sel.TypeText("Reduce ()")

$sel.TypeParagraph()
$sel.TypeText($placeholder)

# Strip the placeholder text out of the blank paragraphs, leaving a clean
# paragraph mark with no run.
$paras = $d.Paragraphs
$count = $paras.Count
For ($i = 1; $i -le $count; $i++) {
  $p = $paras.Item($i)
  $t = $p.Range.Text
  if ($t -eq ($placeholder + "`r")) {
    $r = $p.Range
    $stripped = $d.Range($r.Start, $r.End - 1)
    $stripped.Text = ""
  }
}
